$wb = $excel.ActiveWorkbook

# --- Rename the original sheet to "main" and rework its layout ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "main"

$ws1.Range("B2").Value = "Price"
$ws1.Range("C2").Value = 202
$ws1.Range("B3").Value = "Shares"
$ws1.Range("C3").Value = 50.854999999999997
$ws1.Range("B4").Value = "MC"
$ws1.Range("C4").Formula = "=+C3*C2"
$ws1.Range("B5").Value = "Cash"
$ws1.Range("C5").Formula = "=328+37"
$ws1.Range("B6").Value = "Debt"
$ws1.Range("C6").Formula = "=17+2987"
$ws1.Range("B7").Value = "EV"
$ws1.Range("C7").Formula = "=+C4-C5+C6"

$ws1.Range("C2:C7").NumberFormat = "#,##0"

# Remove the old K:L layout (now unused)
$ws1.Range("K2:L7").Clear()

$ws1.Range("C7").Select()

# --- Add the new "model" sheet after "main" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "model"

# Header row - quarters, filled D..J first, then B,C (matches authoring order)
$ws2.Range("D1").Value = "Q323"
$ws2.Range("E1").Value = "Q423"
$ws2.Range("F1").Value = "Q124"
$ws2.Range("G1").Value = "Q224"
$ws2.Range("H1").Value = "Q324"
$ws2.Range("I1").Value = "Q424"
$ws2.Range("J1").Value = "Q125"
$ws2.Range("B1").Value = "Q123"
$ws2.Range("C1").Value = "Q223"

# Column A line-item labels, in authoring order
$ws2.Range("A3").Value = "Energy and others"
$ws2.Range("A4").Value = "Derivative gains"
$ws2.Range("A5").Value = "Revenue"
$ws2.Range("A6").Value = "Full and energy purchases"
$ws2.Range("A7").Value = "Nuclear fuel amortization "
$ws2.Range("A8").Value = "Derivative losses"
$ws2.Range("A9").Value = "Operating Margins"
$ws2.Range("A11").Value = "GA"
$ws2.Range("A10").Value = "RD"
$ws2.Range("A12").Value = "DA"
$ws2.Range("A13").Value = "Impairments"
$ws2.Range("A14").Value = "Other"
$ws2.Range("A15").Value = "OPEX"
$ws2.Range("A16").Value = "Interest Income "
$ws2.Range("A17").Value = "PRETAX"
$ws2.Range("A18").Value = "TAXES"
$ws2.Range("A19").Value = "NI"
$ws2.Range("A20").Value = "EPS"
$ws2.Range("A21").Value = "Shares"
$ws2.Range("A2").Value = "Capacity rev"

$ws2.Columns.Item(1).ColumnWidth = 24.43

$ws2.Range("I1").Select()
